{"js": "// The buyer (\"Cumparator/Customer\") section listed five detail lines\n// (Denumire, CUI/Tax ID, Adresa, Registrul comertului, Email) for a\n// company buyer. The buyer turned out to be a private individual, so all\n// five lines collapse into a single line: \"Persoana fizica\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Denumire: ...\" paragraph (first line of the buyer block)\n// and the following four paragraphs that make up the rest of the block,\n// by matching on their known text prefixes rather than hard-coded indexes.\nlet denumireIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Denumire:\") === 0) {\n    denumireIndex = i;\n    break;\n  }\n}\n\nif (denumireIndex === -1) {\n  throw new Error('Could not find the \"Denumire:\" paragraph.');\n}\n\n// Replace the first paragraph's text, then delete the next four\n// paragraphs (CUI/Tax ID, Adresa/Adress, Registrul comertului, Email).\nitems[denumireIndex].insertText(\"Persoana fizica\", Word.InsertLocation.replace);\nfor (let i = 1; i <= 4; i++) {\n  items[denumireIndex + i].delete();\n}\n\nawait context.sync();\n", "ps1": "# The buyer (\"Cumparator/Customer\") section listed five detail lines\n# (Denumire, CUI/Tax ID, Adresa, Registrul comertului, Email) for a\n# company buyer. The buyer turned out to be a private individual, so all\n# five lines collapse into a single line: \"Persoana fizica\".\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\n# Locate the \"Denumire: ...\" paragraph (first line of the buyer block).\n$denumireIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    if ($paras.Item($i).Range.Text.StartsWith(\"Denumire:\")) {\n        $denumireIndex = $i\n        break\n    }\n}\n\nif ($denumireIndex -eq -1) {\n    throw \"Could not find the 'Denumire:' paragraph.\"\n}\n\n# Replace that paragraph's text (assigning Range.Text keeps the paragraph\n# mark in place, it just swaps the text before it).\n$paras.Item($denumireIndex).Range.Text = \"Persoana fizica\"\n\n# Delete the next four paragraphs (CUI/Tax ID, Adresa/Adress,\n# Registrul comertului, Email). Each deletion shifts the following\n# paragraph into the same index, so repeatedly deleting the paragraph\n# right after \"Persoana fizica\" removes all four.\nfor ($k = 0; $k -lt 4; $k++) {\n    $d.Paragraphs.Item($denumireIndex + 1).Range.Delete()\n}\n"}
